$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H74").Value = 6399.1816
$ws.Range("I74").Value = 5841.7144
$ws.Range("K74").Value = 5841.7144
$ws.Range("M74").Value = -4905.7144

$ws.Range("H77").Value = 6399.1816
$ws.Range("I77").Value = 5841.7144
$ws.Range("K77").Value = 29208.572
$ws.Range("M77").Value = -24528.572

$ws.Range("H86").Value = 333336640
$ws.Range("I86").Value = 500001000
$ws.Range("K86").Value = 500001000
$ws.Range("M86").Value = -499999877

$ws.Range("H89").Value = 333336640
$ws.Range("I89").Value = 500001000
$ws.Range("K89").Value = 2500005000
$ws.Range("M89").Value = -2499999384

$ws.Range("H100").Value = 3510.3845
$ws.Range("I100").Value = 2411.2
$ws.Range("J100").Value = 4197.375
$ws.Range("K100").Value = 2411.2
$ws.Range("L100").Value = 4197.375
$ws.Range("M100").Value = -1870.2
$ws.Range("N100").Value = -5279.375

$ws.Range("H132").Value = 5601.885
$ws.Range("I132").Value = 5814.5415
$ws.Range("K132").Value = 17443.6245
$ws.Range("M132").Value = -14913.6245

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 5817.1665
$ws.Range("I32").Value = 3479.8333
$ws.Range("K32").Value = 3479.8333
$ws.Range("M32").Value = -3192.8333

$ws.Range("H61").Value = 4229.077
$ws.Range("I61").Value = 3956.5
$ws.Range("J61").Value = 7500
$ws.Range("K61").Value = 3956.5
$ws.Range("L61").Value = 7500
$ws.Range("M61").Value = -3744.5
$ws.Range("N61").Value = -7924

$ws.Range("H97").Value = 1679.8
$ws.Range("J97").Value = 2400.25
$ws.Range("L97").Value = 2400.25
$ws.Range("N97").Value = -3392.25

$ws.Range("H110").Value = 3000.25
$ws.Range("I110").Value = 2750.3333
$ws.Range("K110").Value = 2750.3333
$ws.Range("M110").Value = -705.3332999999998

$ws.Range("H132").Value = 3081.724
$ws.Range("I132").Value = 2709.348
$ws.Range("K132").Value = 8128.044
$ws.Range("M132").Value = -5598.044

$ws.Range("H136").Value = 4229.077
$ws.Range("I136").Value = 3956.5
$ws.Range("J136").Value = 7500
$ws.Range("K136").Value = 11869.5
$ws.Range("L136").Value = 22500
$ws.Range("M136").Value = -9319.5
$ws.Range("N136").Value = -27600

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 1599.6666
$ws.Range("I86").Value = 0
$ws.Range("J86").Value = 1599.6666
$ws.Range("K86").Value = 0
$ws.Range("L86").Value = 1599.6666
$ws.Range("M86").ClearContents()
$ws.Range("N86").Value = -3845.6666

$ws.Range("H89").Value = 1599.6666
$ws.Range("I89").Value = 0
$ws.Range("J89").Value = 1599.6666
$ws.Range("K89").Value = 0
$ws.Range("L89").Value = 7998.333000000001
$ws.Range("M89").ClearContents()
$ws.Range("N89").Value = -19230.333

$ws.Range("H107").Value = 1562.6046
$ws.Range("I107").Value = 1434.3125
$ws.Range("K107").Value = 1434.3125
$ws.Range("M107").Value = 485.6875

$ws.Range("H134").Value = 2664.6428
$ws.Range("I134").Value = 2080.2
$ws.Range("J134").Value = 2989.3333
$ws.Range("K134").Value = 6240.599999999999
$ws.Range("L134").Value = 8967.999899999999
$ws.Range("M134").Value = -3705.599999999999
$ws.Range("N134").Value = -14037.9999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 130.85715
$ws.Range("I7").Value = 67.09090999999999
$ws.Range("J7").Value = 364.66666
$ws.Range("K7").Value = 67.09090999999999
$ws.Range("L7").Value = 364.66666
$ws.Range("M7").Value = 45.90909000000001
$ws.Range("N7").Value = -590.66666

$ws.Range("H31").Value = 3629.4092
$ws.Range("I31").Value = 2128.5557
$ws.Range("K31").Value = 2128.5557
$ws.Range("M31").Value = -1833.5557

$ws.Range("H34").Value = 3629.4092
$ws.Range("I34").Value = 2128.5557
$ws.Range("K34").Value = 2128.5557
$ws.Range("M34").Value = -1926.5557

$ws.Range("H99").Value = 3999.1428
$ws.Range("J99").Value = 3995
$ws.Range("L99").Value = 3995
$ws.Range("N99").Value = -6991

$ws.Range("H107").Value = 85854.836
$ws.Range("J107").Value = 3625
$ws.Range("L107").Value = 3625
$ws.Range("N107").Value = -7465

$ws.Range("H126").Value = 3999.1428
$ws.Range("J126").Value = 3995
$ws.Range("L126").Value = 11985
$ws.Range("N126").Value = -16925

$ws.Range("H134").Value = 3615
$ws.Range("I134").Value = 2908.3125
$ws.Range("J134").Value = 5499.5
$ws.Range("K134").Value = 8724.9375
$ws.Range("L134").Value = 16498.5
$ws.Range("M134").Value = -6189.9375
$ws.Range("N134").Value = -21568.5

$ws.Range("H141").Value = 665589.3
$ws.Range("J141").Value = 665589.3
$ws.Range("L141").Value = 665589.3
$ws.Range("N141").Value = -675949.3

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H128").Value = 115019
$ws.Range("I128").Value = 115019
$ws.Range("K128").Value = 345057
$ws.Range("M128").Value = -340077

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 1361.8182
$ws.Range("J97").Value = 3000
$ws.Range("L97").Value = 3000
$ws.Range("N97").Value = -3992

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 9115.700000000001
$ws.Range("I7").Value = 9363.625
$ws.Range("J7").Value = 8124
$ws.Range("K7").Value = 9363.625
$ws.Range("L7").Value = 8124
$ws.Range("M7").Value = -9251.625
$ws.Range("N7").Value = -8348

$ws.Range("H16").Value = 1849.091
$ws.Range("I16").Value = 1849.091
$ws.Range("K16").Value = 1849.091
$ws.Range("M16").Value = -1679.091

$ws.Range("H46").Value = 2618.111
$ws.Range("I46").Value = 764.2
$ws.Range("J46").Value = 3331.1538
$ws.Range("K46").Value = 764.2
$ws.Range("L46").Value = 3331.1538
$ws.Range("M46").Value = -576.2
$ws.Range("N46").Value = -3707.1538

$ws.Range("H94").Value = 40000
$ws.Range("J94").Value = 40000
$ws.Range("L94").Value = 40000
$ws.Range("N94").Value = -41352

$ws.Range("H112").Value = 131999
$ws.Range("J112").Value = 131999
$ws.Range("L112").Value = 131999
$ws.Range("N112").Value = -134953

$ws.Range("H118").Value = 116000
$ws.Range("J118").Value = 116000
$ws.Range("L118").Value = 116000
$ws.Range("N118").Value = -119314

$ws.Range("H126").Value = 9115.700000000001
$ws.Range("I126").Value = 9363.625
$ws.Range("J126").Value = 8124
$ws.Range("K126").Value = 28090.875
$ws.Range("L126").Value = 24372
$ws.Range("M126").Value = -25620.875
$ws.Range("N126").Value = -29312

$ws.Range("H137").Value = 81496.75
$ws.Range("J137").Value = 81995.664
$ws.Range("L137").Value = 81995.664
$ws.Range("N137").Value = -92195.664

$ws.Range("H139").Value = 78698.25
$ws.Range("J139").Value = 84714.336
$ws.Range("L139").Value = 84714.336
$ws.Range("N139").Value = -94994.336

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 4750
$ws.Range("I62").Value = 3700
$ws.Range("J62").Value = 5800
$ws.Range("K62").Value = 3700
$ws.Range("L62").Value = 5800
$ws.Range("M62").Value = -3076
$ws.Range("N62").Value = -7048

$ws.Range("H65").Value = 4750
$ws.Range("I65").Value = 3700
$ws.Range("J65").Value = 5800
$ws.Range("K65").Value = 18500
$ws.Range("L65").Value = 29000
$ws.Range("M65").Value = -15380
$ws.Range("N65").Value = -35240

$ws.Range("H93").Value = 0
$ws.Range("J93").Value = 0
$ws.Range("L93").Value = 0
$ws.Range("N93").ClearContents()

$ws.Range("H98").Value = 0
$ws.Range("J98").Value = 0
$ws.Range("L98").Value = 0
$ws.Range("N98").ClearContents()

$ws.Range("H99").Value = 0
$ws.Range("J99").Value = 0
$ws.Range("L99").Value = 0
$ws.Range("N99").ClearContents()

$ws.Range("H107").Value = 583.875
$ws.Range("I107").Value = 579.6667
$ws.Range("J107").Value = 586.4
$ws.Range("K107").Value = 1739.0001
$ws.Range("L107").Value = 1759.2
$ws.Range("M107").Value = 180.9999
$ws.Range("N107").Value = -5599.2

$ws.Range("H108").Value = 65000
$ws.Range("J108").Value = 65000
$ws.Range("L108").Value = 65000
$ws.Range("N108").Value = -72680

$ws.Range("H126").Value = 15408.833
$ws.Range("I126").Value = 15408.833
$ws.Range("K126").Value = 46226.499
$ws.Range("M126").Value = -43756.499

$ws.Range("H135").Value = 95999.5
$ws.Range("J135").Value = 95999.5
$ws.Range("L135").Value = 95999.5
$ws.Range("N135").Value = -106139.5

$ws.Range("H137").Value = 55000
$ws.Range("J137").Value = 55000
$ws.Range("L137").Value = 55000
$ws.Range("N137").Value = -65200
